$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-02 Friday" "2024-02-03 Saturday"

Replace-Text "635÷2=317, 1" "881÷7=125, 6"
Replace-Text "239÷4=59, 3" "352÷2=176, 0"
Replace-Text "911÷8=113, 7" "692÷8=86, 4"
Replace-Text "944÷7=134, 6" "939÷8=117, 3"
Replace-Text "948÷4=237, 0" "758÷5=151, 3"

Replace-Text "199÷4=49, 3" "983÷9=109, 2"
Replace-Text "693÷7=99, 0" "511÷7=73, 0"
Replace-Text "680÷9=75, 5" "692÷8=86, 4"
Replace-Text "296÷2=148, 0" "388÷3=129, 1"
Replace-Text "401÷2=200, 1" "364÷9=40, 4"

Replace-Text "569÷4=142, 1" "586÷7=83, 5"
Replace-Text "913÷4=228, 1" "690÷7=98, 4"
Replace-Text "820÷2=410, 0" "598÷5=119, 3"
Replace-Text "756÷5=151, 1" "912÷8=114, 0"
Replace-Text "339÷9=37, 6" "389÷5=77, 4"

Replace-Text "297÷2=148, 1" "716÷3=238, 2"
Replace-Text "177÷8=22, 1" "293÷8=36, 5"
Replace-Text "273÷9=30, 3" "828÷9=92, 0"
Replace-Text "289÷3=96, 1" "846÷7=120, 6"
Replace-Text "921÷8=115, 1" "480÷9=53, 3"

Replace-Text "496÷7=70, 6" "608÷8=76, 0"
Replace-Text "836÷8=104, 4" "269÷8=33, 5"
Replace-Text "694÷7=99, 1" "576÷5=115, 1"
Replace-Text "423÷4=105, 3" "895÷5=179, 0"
Replace-Text "342÷7=48, 6" "298÷5=59, 3"
